$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 4 new product rows (between the existing last data row 18 and the
# totals row, which was row 19 and becomes row 23) and clone row 18's
# formatting (styles + merges) onto them.
# ---------------------------------------------------------------------------
$ws.Range("A19:A22").EntireRow.Insert()

$ws.Range("A18:N18").Copy()
$ws.Range("A19:N22").PasteSpecial(-4122)

$ws.Rows.Item(19).RowHeight = 24.75
$ws.Rows.Item(20).RowHeight = 25.5
$ws.Rows.Item(21).RowHeight = 24.75
$ws.Rows.Item(22).RowHeight = 25.5
$ws.Rows.Item(24).RowHeight = 17.25

$ws.Range("B19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("B20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()
$ws.Range("B21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("B22:G22").Merge()
$ws.Range("H22:K22").Merge()
$ws.Range("L22:M22").Merge()

# ---------------------------------------------------------------------------
# Refresh the full product list (rows 6-22): re-sorted/re-numbered rows
# 6-18 (one new product, BORGASONE, inserted ahead of CETAL; another,
# PULMICORT, inserted ahead of STRINGAZOLE; FUCIDIN's duplicated "1:0" row
# removed) plus four brand-new rows 19-22.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 'BORGASONE TRIO  CREAM '
$ws.Range("H6").Value = '0:0'
$ws.Range("L6").Value = 30.5
$ws.Range("N6").Value = '1:0'
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 'CETAL 250MG/5ML 60ML SUSP'
$ws.Range("H7").Value = '13:0'
$ws.Range("L7").Value = 62
$ws.Range("N7").Value = '2:0'
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 'DOLCYL M 2/500MG 20 F.C. TAB'
$ws.Range("H8").Value = '0:1'
$ws.Range("L8").Value = 17
$ws.Range("N8").Value = '0:2'
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 'FUCIDIN 2% CREAM 30 GM'
$ws.Range("H9").Value = '1:0'
$ws.Range("L9").Value = -96
$ws.Range("N9").Value = '1:0'
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = 'GLYBOFEN 5/500MG 30 F.C.TABS.'
$ws.Range("H10").Value = '1:2'
$ws.Range("L10").Value = -15
$ws.Range("N10").Value = '0:0'
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 'HIBIOTIC 1GM 16 TAB'
$ws.Range("H11").Value = '2:0'
$ws.Range("L11").Value = 86.5
$ws.Range("N11").Value = '0:2'
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 'KETOLAC 10MG 20 TAB'
$ws.Range("H12").Value = '0:1'
$ws.Range("L12").Value = 19
$ws.Range("N12").Value = '0:2'
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 'OPLEX-N SYRUP 125ML'
$ws.Range("H13").Value = '6:0'
$ws.Range("L13").Value = 31
$ws.Range("N13").Value = '1:0'
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = 'OTAL EAR DROPS 5 ML'
$ws.Range("H14").Value = '3:0'
$ws.Range("L14").Value = 19
$ws.Range("N14").Value = '1:0'
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = 'PANADOL COLD & FLU DAY 24 F.C. TABS'
$ws.Range("H15").Value = '2:0'
$ws.Range("L15").Value = 76
$ws.Range("N15").Value = '1:0'
$ws.Range("A16").Value = 13
$ws.Range("B16").Value = 'PULMICORT 0.25MG/ML 20 NEBULIZER VIAL SUSP.'
$ws.Range("H16").Value = '0:7'
$ws.Range("L16").Value = 56.4
$ws.Range("N16").Value = '0:0'
$ws.Range("A17").Value = 14
$ws.Range("B17").Value = 'STRINGAZOLE 40MG 21 ENTERIC COATED TABLETS'
$ws.Range("H17").Value = '1:0'
$ws.Range("L17").Value = 42
$ws.Range("N17").Value = '0:0'
$ws.Range("A18").Value = 15
$ws.Range("B18").Value = 'URICONTROL 2MG 10 TAB.'
$ws.Range("H18").Value = '1:0'
$ws.Range("L18").Value = 25
$ws.Range("N18").Value = '1:0'
$ws.Range("A19").Value = 16
$ws.Range("B19").Value = 'سرنجات 3 سم'
$ws.Range("H19").Value = '-1:0'
$ws.Range("L19").Value = 2
$ws.Range("N19").Value = '1:0'
$ws.Range("A20").Value = 17
$ws.Range("B20").Value = 'كريم شعر نيو هير صغير'
$ws.Range("H20").Value = '2:0'
$ws.Range("L20").Value = 20
$ws.Range("N20").Value = '1:0'
$ws.Range("A21").Value = 18
$ws.Range("B21").Value = 'لزقه النمر بسعر القطعه'
$ws.Range("H21").Value = '48:0'
$ws.Range("L21").Value = 15
$ws.Range("N21").Value = '1:0'
$ws.Range("A22").Value = 19
$ws.Range("B22").Value = 'مبرد قدم'
$ws.Range("H22").Value = '1:0'
$ws.Range("L22").Value = 40
$ws.Range("N22").Value = '2:0'

# ---------------------------------------------------------------------------
# Recompute the grand total (K column) now that the product list changed.
# ---------------------------------------------------------------------------
$ws.Range("K23").Value = 486.4
